$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.843.19"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "3.162.12"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "3.132.74"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("E9").Value = "  -2.36%  "

$ws.Range("E10").Value = "  -4.88%  "

$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("E12").Value = "  -3.21%  "

$ws.Range("E13").Value = "  -5.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.15%  "

$ws.Range("D15").Value = "3.683.70"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.148.13"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.659.42"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("E23").Value = "  -2.23%  "

$ws.Range("E24").Value = "  -4.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.54%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  -4.19%  "

$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("E37").Value = "  -4.96%  "

$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "434.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.42%  "

$ws.Range("D43").Value = "2.911.21"
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("E44").Value = "  -2.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.280"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("E46").Value = "  -4.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("E50").Value = "  -2.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.67%  "

